# Update crypto price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "33.713.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.774.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.558"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.30%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.23"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.61"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.58%  "
$ws.Range("E10").Value = "  +2.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0662"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.00%  "
$ws.Range("E12").Value = "  +1.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.031.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.779.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.622"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "33.714.87"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "9.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.99%  "
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.41"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "250.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0736"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.23"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("E24").Value = "  -2.27%  "
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.25%  "
$ws.Range("E28").Value = "  +0.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("E32").Value = "  +2.02%  "
$ws.Range("E33").Value = "  +3.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.53%  "
$ws.Range("E35").Value = "  +4.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.475.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.07"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.626"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0184"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.24%  "
$ws.Range("E41").Value = "  +1.79%  "
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.881"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.928.31"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.26%  "
$ws.Range("E48").Value = "  +0.16%  "
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.65%  "
